$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7192
$ws.Range("C3").Value = 171634
$ws.Range("C4").Value = 162436
$ws.Range("C8").Value = 65.87
